$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts existing D:K quarterly data to E:L)
# to make room for the newest reported quarter (period ending 2018-09-30, serial 43373).
$ws.Columns("D").Insert()

# Copy number formats / styles from the (now-shifted) old column D, which is
# now column E, into the freshly inserted column D, per contiguous data block.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column D with the newest quarter's figures.

# Income Statement block
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 51500
$ws.Range("D9").Value = 18200
$ws.Range("D10").Value = 33300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 200
$ws.Range("D17").Value = 67500
$ws.Range("D18").Value = -16000
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -15800
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -16000
$ws.Range("D24").Value = 100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -16200
$ws.Range("D27").Value = -16100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -16100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -16100

# Balance Sheet block
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 29300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 16100
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 382500
$ws.Range("D48").Value = 2700
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 30400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 461000
$ws.Range("D57").Value = 107300
$ws.Range("D58").Value = 293300
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 108200
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 508500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -261800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -47500
$ws.Range("D77").Value = 0

# Cash Flow Statement block
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -16100
$ws.Range("D83").Value = 200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -12300
$ws.Range("D91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -30700
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 23800
$ws.Range("D101").Value = 200
$ws.Range("D102").Value = -19100
